# Update "想去人数" (F column) counts on each sheet to reflect the
# regenerated gh-pages data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 7927  # F5: 7900 -> 7927
$ws.Cells.Item(8, 6).Value = 80  # F8: 81 -> 80
$ws.Cells.Item(9, 6).Value = 52  # F9: 47 -> 52
$ws.Cells.Item(10, 6).Value = 6837  # F10: 6804 -> 6837
$ws.Cells.Item(11, 6).Value = 1126  # F11: 1123 -> 1126
$ws.Cells.Item(12, 6).Value = 499  # F12: 494 -> 499
$ws.Cells.Item(13, 6).Value = 471  # F13: 464 -> 471
$ws.Cells.Item(15, 6).Value = 677  # F15: 670 -> 677
$ws.Cells.Item(16, 6).Value = 343  # F16: 342 -> 343
$ws.Cells.Item(17, 6).Value = 301  # F17: 300 -> 301
$ws.Cells.Item(19, 6).Value = 74  # F19: 73 -> 74
$ws.Cells.Item(20, 6).Value = 165  # F20: 164 -> 165
$ws.Cells.Item(22, 6).Value = 11134  # F22: 11067 -> 11134
$ws.Cells.Item(23, 6).Value = 97  # F23: 96 -> 97
$ws.Cells.Item(24, 6).Value = 79  # F24: 71 -> 79
$ws.Cells.Item(25, 6).Value = 2122  # F25: 2104 -> 2122
$ws.Cells.Item(26, 6).Value = 2838  # F26: 2785 -> 2838
$ws.Cells.Item(28, 6).Value = 42  # F28: 41 -> 42
$ws.Cells.Item(29, 6).Value = 2515  # F29: 2491 -> 2515
$ws.Cells.Item(32, 6).Value = 36  # F32: 35 -> 36
$ws.Cells.Item(34, 6).Value = 2280  # F34: 2257 -> 2280
$ws.Cells.Item(36, 6).Value = 1554  # F36: 1535 -> 1554
$ws.Cells.Item(37, 6).Value = 68  # F37: 66 -> 68
$ws.Cells.Item(38, 6).Value = 64  # F38: 56 -> 64
$ws.Cells.Item(39, 6).Value = 5628  # F39: 5594 -> 5628
$ws.Cells.Item(40, 6).Value = 69  # F40: 66 -> 69
$ws.Cells.Item(41, 6).Value = 1236  # F41: 1231 -> 1236
$ws.Cells.Item(42, 6).Value = 801  # F42: 797 -> 801
$ws.Cells.Item(43, 6).Value = 149  # F43: 147 -> 149
$ws.Cells.Item(45, 6).Value = 1094  # F45: 1095 -> 1094
$ws.Cells.Item(46, 6).Value = 1045  # F46: 1043 -> 1045
$ws.Cells.Item(47, 6).Value = 1470  # F47: 1462 -> 1470
$ws.Cells.Item(48, 6).Value = 90  # F48: 84 -> 90
$ws.Cells.Item(49, 6).Value = 1118  # F49: 1117 -> 1118

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(8, 6).Value = 247  # F8: 246 -> 247
$ws.Cells.Item(9, 6).Value = 9  # F9: 8 -> 9
$ws.Cells.Item(10, 6).Value = 51  # F10: 50 -> 51
$ws.Cells.Item(11, 6).Value = 207  # F11: 206 -> 207
$ws.Cells.Item(20, 6).Value = 55  # F20: 52 -> 55
$ws.Cells.Item(23, 6).Value = 4  # F23: 3 -> 4

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 154  # F2: 141 -> 154
$ws.Cells.Item(3, 6).Value = 261  # F3: 245 -> 261

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 154  # F4: 141 -> 154
$ws.Cells.Item(5, 6).Value = 261  # F5: 245 -> 261
$ws.Cells.Item(8, 6).Value = 7927  # F8: 7900 -> 7927
$ws.Cells.Item(9, 6).Value = 80  # F9: 81 -> 80
$ws.Cells.Item(11, 6).Value = 52  # F11: 47 -> 52
$ws.Cells.Item(12, 6).Value = 6838  # F12: 6804 -> 6838
$ws.Cells.Item(13, 6).Value = 6838  # F13: 6804 -> 6838
$ws.Cells.Item(14, 6).Value = 1126  # F14: 1123 -> 1126
$ws.Cells.Item(15, 6).Value = 499  # F15: 494 -> 499
$ws.Cells.Item(16, 6).Value = 471  # F16: 464 -> 471
$ws.Cells.Item(17, 6).Value = 677  # F17: 670 -> 677
$ws.Cells.Item(18, 6).Value = 343  # F18: 342 -> 343
$ws.Cells.Item(19, 6).Value = 301  # F19: 300 -> 301
$ws.Cells.Item(21, 6).Value = 247  # F21: 246 -> 247
$ws.Cells.Item(22, 6).Value = 165  # F22: 164 -> 165
$ws.Cells.Item(23, 6).Value = 207  # F23: 206 -> 207
$ws.Cells.Item(25, 6).Value = 11134  # F25: 11067 -> 11134
$ws.Cells.Item(26, 6).Value = 97  # F26: 96 -> 97
$ws.Cells.Item(27, 6).Value = 79  # F27: 71 -> 79
$ws.Cells.Item(28, 6).Value = 2122  # F28: 2104 -> 2122
$ws.Cells.Item(29, 6).Value = 2838  # F29: 2785 -> 2838
$ws.Cells.Item(30, 6).Value = 2515  # F30: 2491 -> 2515
$ws.Cells.Item(33, 6).Value = 36  # F33: 35 -> 36
$ws.Cells.Item(35, 6).Value = 2280  # F35: 2257 -> 2280
$ws.Cells.Item(37, 6).Value = 1554  # F37: 1536 -> 1554
$ws.Cells.Item(38, 6).Value = 68  # F38: 66 -> 68
$ws.Cells.Item(39, 6).Value = 64  # F39: 56 -> 64
$ws.Cells.Item(40, 6).Value = 5628  # F40: 5594 -> 5628
$ws.Cells.Item(41, 6).Value = 55  # F41: 52 -> 55
$ws.Cells.Item(42, 6).Value = 1236  # F42: 1231 -> 1236
$ws.Cells.Item(43, 6).Value = 801  # F43: 797 -> 801
$ws.Cells.Item(44, 6).Value = 149  # F44: 147 -> 149
$ws.Cells.Item(46, 6).Value = 1094  # F46: 1095 -> 1094
$ws.Cells.Item(47, 6).Value = 1045  # F47: 1043 -> 1045
$ws.Cells.Item(48, 6).Value = 1470  # F48: 1462 -> 1470
$ws.Cells.Item(49, 6).Value = 90  # F49: 84 -> 90
$ws.Cells.Item(50, 6).Value = 1118  # F50: 1117 -> 1118

